$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update translated Telugu text for records that were re-worded ---
$ws.Range("D5").Value = "అరంగేట్రం చేసిన పురాతన ( ఓల్డ్ ) ఆటగాళ్ల"
$ws.Range("D19").Value = "కెప్టెన్ గా ఒక ఇన్నింగ్స్ లో ఉత్తమ గణాంకాలు సాధించిన ఆటగాళ్ల"
$ws.Range("D26").Value = "చివరి మ్యాచ్ లో శతకం చేసిన ఆటగాళ్ల"
$ws.Range("D30").Value = "వరుసగా అత్యధిక డక్లు సాధించిన ఆటగాళ్ల"
$ws.Range("D39").Value = "కెరీర్ లో అత్యధిక మెయిడెన్ లు వేసిన ఆటగాళ్ల"
$ws.Range("D43").Value = "కెరీర్ లో అతి తక్కువ సార్లు డక్ అవుట్ అయిన ఆటగాళ్ల"
$ws.Range("D49").Value = "ఒక ఇన్నింగ్స్ లో అత్యధిక మెయిడెన్ లు వేసిన ఆటగాళ్ల"
$ws.Range("D52").Value = "అత్యధిక ప్లేయర్-ఆఫ్-ది-సిరీస్ అవార్డులు గెల్చుకున్న ఆటగాళ్ల"
$ws.Range("D53").Value = "ఒక ఇన్నింగ్స్ లో ఐదు వికెట్లు సాధించిన పిన్న వయసు ఆటగాళ్ల"
$ws.Range("D56").Value = "మొదటి మ్యాచ్ లో సెంచరీ సాధించిన పురాతన ( ఓల్డ్ ) ఆటగాళ్ల"
$ws.Range("D64").Value = "డక్ అవుట్ అవ్వకుండ అత్యధిక ఇన్నింగ్స్ ఆడిన ఆటగాళ్ల"
$ws.Range("D68").Value = "పిన్న వయసులో కెప్టెన్ గా వ్యావహరించిన ఆటగాళ్ల"
$ws.Range("D69").Value = "ఒక జట్టుకి  కెప్టెన్ గా అత్యధిక వరుస మ్యాచ్ లలో ఆడిన ఆటగాళ్ల"

# --- Restore default font (Calibri) on cells that previously had the blank/Arial font ---
$calibriCells = @(
    "D3",
    "D7",
    "D10",
    "E10",
    "D11",
    "E12",
    "D13",
    "D14",
    "E14",
    "D16",
    "D17",
    "E17",
    "D18",
    "E18",
    "D19",
    "D21",
    "D24",
    "E24",
    "D26",
    "D27",
    "E27",
    "E30",
    "D31",
    "E31",
    "D32",
    "D34",
    "D35",
    "E37",
    "E40",
    "D41",
    "D42",
    "E43",
    "D50",
    "D54",
    "E54",
    "D55",
    "E55",
    "E56",
    "E59",
    "E64",
    "D67",
    "E67",
    "E69",
    "F69",
    "D71",
    "D72",
    "E72",
    "D73",
    "E76",
    "D78",
    "D79",
    "E79",
    "E81"
)
foreach ($addr in $calibriCells) {
    $ws.Range($addr).Font.Name = "Calibri"
}

# --- Reset a few cells (whose rows were re-translated) back to the plain/default font ---
$arialCells = @(
    "D52",
    "D53",
    "D56",
    "D68"
)
foreach ($addr in $arialCells) {
    $ws.Range($addr).Font.Name = "Arial"
}
